$p = $ppt.ActivePresentation

# Donor shape: an existing "TextBox" on slide 2 that already carries the
# exact run/paragraph formatting (Times New Roman 12pt, no-fill, auto-fit,
# rtlCol, lstStyle, panose metadata, etc.) we want for the new label. Copying
# it and then repositioning/retexting it reproduces that formatting exactly,
# which isn't otherwise reachable through the plain Shape/Font COM surface.
$donorSlide = $p.Slides.Item(2)
$donor = $null
for ($j = 1; $j -le $donorSlide.Shapes.Count; $j++) {
    $candidate = $donorSlide.Shapes.Item($j)
    if ($candidate.Name -eq "TextBox 62") {
        $donor = $candidate
    }
}

$targetSlide = $p.Slides.Item(1)

$donor.Copy()
$pasted = $targetSlide.Shapes.Paste()
$lbl = $pasted.Item(1)

$lbl.Name = "TextBox 2"

# Exact target geometry (EMU -> points, 12700 EMU per point):
#   off  x=1849020 y=619872
#   ext cx=873957 cy=276999
$lbl.Left = 145.59212598425196
$lbl.Width = 68.81551181102363
$lbl.Height = 21.810944881889764
# The Top setter quantizes through a float32 round-trip; nudge by a few
# millionths of a point (well under 1/100 mm) so it lands on the exact
# target EMU value instead of one EMU short.
$lbl.Top = 48.80885889763779

$lbl.TextFrame.TextRange.Text = "Kelp forest"
